$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 14-35 must be stored as text (matching the source t="inlineStr" cells),
# so force a text number format before assignment, then clear the format override
# afterward to avoid leaving a stray style index on cells that should use the default style.
$newDataRange = $ws.Range("A14:E35")
$newDataRange.NumberFormat = "@"

# Row 14: Can - Olive Black (Sliced)
$ws.Cells.Item(14, 1).Value = "62063"
$ws.Cells.Item(14, 2).Value = "Can - Olive Black (Sliced)"
$ws.Cells.Item(14, 3).Value = "1"
$ws.Cells.Item(14, 4).Value = "$47.95"
$ws.Cells.Item(14, 5).Value = "$47.95"

# Row 15: Green Bean - FRZ Whole
$ws.Cells.Item(15, 1).Value = "45340"
$ws.Cells.Item(15, 2).Value = "Green Bean - FRZ Whole"
$ws.Cells.Item(15, 3).Value = "1"
$ws.Cells.Item(15, 4).Value = "$32.05"
$ws.Cells.Item(15, 5).Value = "$32.05"

# Row 16: Broccoli - FRZ
$ws.Cells.Item(16, 1).Value = "1059229"
$ws.Cells.Item(16, 2).Value = "Broccoli - FRZ"
$ws.Cells.Item(16, 3).Value = "1"
$ws.Cells.Item(16, 4).Value = "$32.78"
$ws.Cells.Item(16, 5).Value = "$32.78"

# Row 17: Cheddar - (Sliced)
$ws.Cells.Item(17, 1).Value = "25806"
$ws.Cells.Item(17, 2).Value = "Cheddar - (Sliced)"
$ws.Cells.Item(17, 3).Value = "15"
$ws.Cells.Item(17, 4).Value = "$33.76"
$ws.Cells.Item(17, 5).Value = "$506.40"

# Row 18: Feta - Pail
$ws.Cells.Item(18, 1).Value = "29053"
$ws.Cells.Item(18, 2).Value = "Feta - Pail"
$ws.Cells.Item(18, 3).Value = "1"
$ws.Cells.Item(18, 4).Value = "$95.45"
$ws.Cells.Item(18, 5).Value = "$95.45"

# Row 19: Lemon Juice
$ws.Cells.Item(19, 1).Value = "46999"
$ws.Cells.Item(19, 2).Value = "Lemon Juice"
$ws.Cells.Item(19, 3).Value = "1"
$ws.Cells.Item(19, 4).Value = "$16.55"
$ws.Cells.Item(19, 5).Value = "$16.55"

# Row 20: Parmesan (Grated)
$ws.Cells.Item(20, 1).Value = "25330"
$ws.Cells.Item(20, 2).Value = "Parmesan (Grated)"
$ws.Cells.Item(20, 3).Value = "2"
$ws.Cells.Item(20, 4).Value = "$59.54"
$ws.Cells.Item(20, 5).Value = "$119.08"

# Row 21: Pickle - Dill Chip
$ws.Cells.Item(21, 1).Value = "60171"
$ws.Cells.Item(21, 2).Value = "Pickle - Dill Chip"
$ws.Cells.Item(21, 3).Value = "1"
$ws.Cells.Item(21, 4).Value = "$34.09"
$ws.Cells.Item(21, 5).Value = "$34.09"

# Row 22: Sauerkraut
$ws.Cells.Item(22, 1).Value = "3275539"
$ws.Cells.Item(22, 2).Value = "Sauerkraut"
$ws.Cells.Item(22, 3).Value = "2"
$ws.Cells.Item(22, 4).Value = "$19.35"
$ws.Cells.Item(22, 5).Value = "$38.70"

# Row 23: Sausage - Chicken Patty
$ws.Cells.Item(23, 1).Value = "2825368"
$ws.Cells.Item(23, 2).Value = "Sausage - Chicken Patty"
$ws.Cells.Item(23, 3).Value = "2"
$ws.Cells.Item(23, 4).Value = "$51.14"
$ws.Cells.Item(23, 5).Value = "$102.28"

# Row 24: Smoked Turkey (Unsliced)
$ws.Cells.Item(24, 1).Value = "54112"
$ws.Cells.Item(24, 2).Value = "Smoked Turkey (Unsliced)"
$ws.Cells.Item(24, 3).Value = "4"
$ws.Cells.Item(24, 4).Value = "$76.44"
$ws.Cells.Item(24, 5).Value = "$305.76"

# Row 25: Sour Cream
$ws.Cells.Item(25, 1).Value = "1132582"
$ws.Cells.Item(25, 2).Value = "Sour Cream"
$ws.Cells.Item(25, 3).Value = "1"
$ws.Cells.Item(25, 4).Value = "$28.94"
$ws.Cells.Item(25, 5).Value = "$28.94"

# Row 26: Tuna White Chunk (Pouch)
$ws.Cells.Item(26, 1).Value = "8255796"
$ws.Cells.Item(26, 2).Value = "Tuna White Chunk (Pouch)"
$ws.Cells.Item(26, 3).Value = "6"
$ws.Cells.Item(26, 4).Value = "$72.73"
$ws.Cells.Item(26, 5).Value = "$436.38"

# Row 27: Vegan Chicken Tenders
$ws.Cells.Item(27, 1).Value = "11072"
$ws.Cells.Item(27, 2).Value = "Vegan Chicken Tenders"
$ws.Cells.Item(27, 3).Value = "2"
$ws.Cells.Item(27, 4).Value = "$87.80"
$ws.Cells.Item(27, 5).Value = "$175.60"

# Row 28: Yogurt - Greek (Bulk)
$ws.Cells.Item(28, 1).Value = "6364494"
$ws.Cells.Item(28, 2).Value = "Yogurt - Greek (Bulk)"
$ws.Cells.Item(28, 3).Value = "2"
$ws.Cells.Item(28, 4).Value = "$27.11"
$ws.Cells.Item(28, 5).Value = "$54.22"

# Row 29: BBQ - Sauce
$ws.Cells.Item(29, 1).Value = "1030192"
$ws.Cells.Item(29, 2).Value = "BBQ - Sauce"
$ws.Cells.Item(29, 3).Value = "1"
$ws.Cells.Item(29, 4).Value = "$72.94"
$ws.Cells.Item(29, 5).Value = "$72.94"

# Row 30: Bacon (Pre-Cooked)
$ws.Cells.Item(30, 1).Value = "5514021"
$ws.Cells.Item(30, 2).Value = "Bacon (Pre-Cooked)"
$ws.Cells.Item(30, 3).Value = "15"
$ws.Cells.Item(30, 4).Value = "$33.50"
$ws.Cells.Item(30, 5).Value = "$502.50"

# Row 31: Arugula - Fresh
$ws.Cells.Item(31, 1).Value = "88202"
$ws.Cells.Item(31, 2).Value = "Arugula - Fresh"
$ws.Cells.Item(31, 3).Value = "4"
$ws.Cells.Item(31, 4).Value = "$18.42"
$ws.Cells.Item(31, 5).Value = "$73.68"

# Row 32: Carrots - Jumbo Fresh
$ws.Cells.Item(32, 1).Value = "7228448"
$ws.Cells.Item(32, 2).Value = "Carrots - Jumbo Fresh"
$ws.Cells.Item(32, 3).Value = "1"
$ws.Cells.Item(32, 4).Value = "$25.10"
$ws.Cells.Item(32, 5).Value = "$25.10"

# Row 33: Cilantro - Fresh
$ws.Cells.Item(33, 1).Value = "6579288"
$ws.Cells.Item(33, 2).Value = "Cilantro - Fresh"
$ws.Cells.Item(33, 3).Value = "1"
$ws.Cells.Item(33, 4).Value = "$24.70"
$ws.Cells.Item(33, 5).Value = "$24.70"

# Row 34: Garlic - Fresh (Peeled)
$ws.Cells.Item(34, 1).Value = "5365192"
$ws.Cells.Item(34, 2).Value = "Garlic - Fresh (Peeled)"
$ws.Cells.Item(34, 3).Value = "1"
$ws.Cells.Item(34, 4).Value = "$68.97"
$ws.Cells.Item(34, 5).Value = "$68.97"

# Row 35: Onion - Yellow Fresh
$ws.Cells.Item(35, 1).Value = "88153"
$ws.Cells.Item(35, 2).Value = "Onion - Yellow Fresh"
$ws.Cells.Item(35, 3).Value = "2"
$ws.Cells.Item(35, 4).Value = "$19.72"
$ws.Cells.Item(35, 5).Value = "$39.44"

# Remove the temporary number-format override now that the text values are locked in,
# restoring the default (General) style on these cells.
$newDataRange.ClearFormats()
